# Weekly update: insert a new price record for "Cilantro" at Vega Modelo de
# Temuco (row 294), shifting the existing rows 294-310 down to 295-311.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 294 (shifts 294..310 -> 295..311)
$ws.Rows.Item(294).Insert()

# Populate the newly inserted row 294 with the new weekly record.
# Category/location columns mirror the surrounding rows for this subset;
# only the date (D) and volume (J) differ for this new entry.
$ws.Cells.Item(294, 1).Value = 10
$ws.Cells.Item(294, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(294, 3).Value = "La Araucanía"
$ws.Cells.Item(294, 4).Value = 44610
$ws.Cells.Item(294, 5).Value = 9
$ws.Cells.Item(294, 6).Value = 100112040
$ws.Cells.Item(294, 7).Value = "Cilantro"
$ws.Cells.Item(294, 8).Value = "Sin especificar"
$ws.Cells.Item(294, 9).Value = "Primera"
$ws.Cells.Item(294, 10).Value = 40
$ws.Cells.Item(294, 11).Value = 5000
$ws.Cells.Item(294, 12).Value = 5000
$ws.Cells.Item(294, 13).Value = 5000
$ws.Cells.Item(294, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(294, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(294, 16).Value = 2500
$ws.Cells.Item(294, 17).Value = 2
$ws.Cells.Item(294, 18).Value = "Hortaliza"
